# Generate Report for Handoff
# Regenerates the localization-status report: new handoff UUID/filenames
# and refreshed handoff timestamps across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "095c778e-6bcd-42e9-ac77-1657a8dfa787"
$newGuid = "fa92cefb-154d-4064-bac7-829c38c5261b"

$oldHash = "c267fc237ee45cf264bf07a931fe1c869e702282"
$newHash = "9e698e62a06d898a4486435e6a956ad5b14a0a67"

$oldMdName  = "$oldGuid.md"
$newMdName  = "$newGuid.md"

$oldZhName  = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhName  = "$newGuid.$newHash.zh-cn.xlf"

$oldDeName  = "$oldGuid.$oldHash.de-de.xlf"
$newDeName  = "$newGuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = "2016-40-11 16:40:46"

$overviewUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0eca65f26228b3d2cb5d1e3b74d809fc1c836287/e2e/$newMdName"
$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $overviewUrl, "", "", $newMdName)
$wsOverview.Range("A2").Font.Color = 15570276
$wsOverview.Range("A2").Font.Underline = 2

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhName
$wsZh.Range("E2").Value = "2016-03-11 16:40:42"

$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0eca65f26228b3d2cb5d1e3b74d809fc1c836287/e2e/$newMdName"
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhMdUrl, "", "", $newMdName)
$wsZh.Range("A2").Font.Color = 15570276
$wsZh.Range("A2").Font.Underline = 2

$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2ec292765cb735104ef3b5016bf6df3617ead6cd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhName"
$wsZh.Range("D2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfUrl, "", "", $newZhName)
$wsZh.Range("D2").Font.Color = 15570276
$wsZh.Range("D2").Font.Underline = 2

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeName
$wsDe.Range("E2").Value = "2016-03-11 16:40:46"

$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0eca65f26228b3d2cb5d1e3b74d809fc1c836287/e2e/$newMdName"
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deMdUrl, "", "", $newMdName)
$wsDe.Range("A2").Font.Color = 15570276
$wsDe.Range("A2").Font.Underline = 2

$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0548d97587dc9e88dfba6a2d8328b31736ddfbdb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeName"
$wsDe.Range("D2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfUrl, "", "", $newDeName)
$wsDe.Range("D2").Font.Color = 15570276
$wsDe.Range("D2").Font.Underline = 2
